# Fruta / hortaliza, semanal
# A new weekly price record (row 120, date 2023-08-03 / serial 45141) is
# inserted into the "Vega Monumental Concepción - Mango" sheet. All the
# rows that used to be 120..184 shift down by one to 121..185.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 120; Excel shifts rows 120:184 down
# to 121:185 and extends the used range to A1:T185.
$ws.Rows("120:120").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A120").Value = 11
$ws.Range("B120").Value = "Vega Monumental Concepción"
$ws.Range("C120").Value = "Bíobío"
$ws.Range("D120").Value = 45141
$ws.Range("E120").Value = 8
$ws.Range("F120").Value = "Fruta"
$ws.Range("G120").Value = 100108
$ws.Range("H120").Value = "Tropicales y subtropicales"
$ws.Range("I120").Value = 100108002
$ws.Range("J120").Value = "Mango"
$ws.Range("K120").Value = "Sin especificar"
$ws.Range("L120").Value = "Primera"
$ws.Range("M120").Value = 100
$ws.Range("N120").Value = 8000
$ws.Range("O120").Value = 8500
$ws.Range("P120").Value = 8250
$ws.Range("Q120").Value = "`$/bandeja 4 kilos"
$ws.Range("R120").Value = "Brasil"
$ws.Range("S120").Value = 2062
$ws.Range("T120").Value = 4
